$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.011.86'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.886.70'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.36%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('E5').Value = '  -2.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4592'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -2.97%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4059'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.86'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.20%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07969'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9918'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.21%  '
$ws.Range('E12').Value = '  -3.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.883.56'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.905'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -3.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.072'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -4.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.30'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -3.66%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001029'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.39%  '
$ws.Range('E19').Value = '  -1.16%  '
$ws.Range('E20').Value = '  -2.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.0000'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '29.010.09'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.414'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.76%  '
$ws.Range('E24').Value = '  +1.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.206'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.083.47'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.52'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.57'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.102'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.423'
$ws.Range('D30').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.79'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.002'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.64%  '
$ws.Range('E33').Value = '  -2.65%  '
$ws.Range('E34').Value = '  -1.55%  '
$ws.Range('E35').Value = '  -0.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.278'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.78%  '
$ws.Range('E37').Value = '  -2.54%  '
$ws.Range('E38').Value = '  -3.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.252'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.80%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.176'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9993'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5781'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.47%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1821'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.43%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.12'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.258'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.95%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.07514'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.88%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.263'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +5.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '12.04'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5453'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.56%  '
$ws.Range('E50').Value = '  -4.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '111.15'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.86%  '
